$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.300.17'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.68%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.097.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.04%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '524.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.16%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.55%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.095.52'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.12%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.446'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.90%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.28'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.09%  '

# Row 11
$ws.Range('E11').Value = '  -1.04%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.393'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.28%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.643.96'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.25%  '

# Row 14
$ws.Range('E14').Value = '  +2.80%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.04%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000162'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.87%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.417.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.63%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.108.26'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.15%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.91'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.78%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.34'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.35%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.45%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '344.33'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.74%  '

# Row 23
$ws.Range('E23').Value = '  -0.03%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.62'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.67%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.498'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.43%  '

# Row 26
$ws.Range('E26').Value = '  -1.64%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.43%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0887'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.95%  '

# Row 29
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'USDe'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.05%  '

# Row 30
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.37'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.93%  '

# Row 31
$ws.Range('E31').Value = '  -0.07%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.85%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.77'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.64%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.91'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +6.84%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.14'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.85%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '158.13'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.44%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.19%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.81'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.33%  '

# Row 39
$ws.Range('E39').Value = '  -1.94%  '

# Row 40
$ws.Range('E40').Value = '  +5.52%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0657'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.40%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.08'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.31%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.697'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.23%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.147.66'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.23%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.374.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.50%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '36.52'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.72%  '

# Row 47
$ws.Range('E47').Value = '  +0.07%  '

# Row 48
$ws.Range('E48').Value = '  +3.09%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.969'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.47%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.96'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.94%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.65'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.95%  '
